# Add the "ODI Bowling Extra" worksheet (scraped extra bowling attributes),
# placed after the existing "ODI Batting Extra" sheet, as the new 5th tab.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Data values are textual in the source export (MATCH_CODE, MAIDEN_OVERS and
# PERCENT_WICKETS_OF_ALL are all stored as text, not numbers/percentages),
# so force the used range to Text format before writing any values.
$usedRange = $ws.Range("A1:C21")
$usedRange.NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "MAIDEN_OVERS"
$ws.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data rows
$ws.Range("A2").Value = "4117"
$ws.Range("B2").Value = "0"

$ws.Range("A3").Value = "4123"
$ws.Range("B3").Value = "0"

$ws.Range("A4").Value = "4137"

$ws.Range("A5").Value = "4139"

$ws.Range("A6").Value = "4149"
$ws.Range("B6").Value = "0"

$ws.Range("A7").Value = "4166"
$ws.Range("B7").Value = "0"

$ws.Range("A8").Value = "4167"

$ws.Range("A9").Value = "4168"
$ws.Range("B9").Value = "0"

$ws.Range("A10").Value = "4169"

$ws.Range("A11").Value = "4170"
$ws.Range("B11").Value = "0"

$ws.Range("A12").Value = "4171"
$ws.Range("B12").Value = "0"

$ws.Range("A13").Value = "4287"

$ws.Range("A14").Value = "4321"
$ws.Range("B14").Value = "0"
$ws.Range("C14").Value = "20.00%"

$ws.Range("A15").Value = "4331"
$ws.Range("B15").Value = "0"

$ws.Range("A16").Value = "4346"

$ws.Range("A17").Value = "4401"

$ws.Range("A18").Value = "4405"

$ws.Range("A19").Value = "4408"

$ws.Range("A20").Value = "4431"
$ws.Range("B20").Value = "0"
$ws.Range("C20").Value = "20.00%"

$ws.Range("A21").Value = "4618"
$ws.Range("B21").Value = "0"
